$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "61.820.10"
$ws.Range("E2").Value = "  +8.49%  "

# Row 3
$ws.Range("D3").Value = "3.421.08"
$ws.Range("E3").Value = "  +5.56%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").Value = "'418.47"
$ws.Range("E5").Value = "  +5.88%  "

# Row 6
$ws.Range("D6").Value = "'117.91"
$ws.Range("E6").Value = "  +8.80%  "

# Row 7
$ws.Range("D7").Value = "3.408.51"
$ws.Range("E7").Value = "  +5.31%  "

# Row 8
$ws.Range("D8").Value = "'0.572"
$ws.Range("E8").Value = "  -2.02%  "

# Row 9
$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "  +0.06%  "

# Row 10
$ws.Range("D10").Value = "'0.627"
$ws.Range("E10").Value = "  -0.17%  "

# Row 11
$ws.Range("D11").Value = "'0.114"
$ws.Range("E11").Value = "  +17.55%  "

# Row 12
$ws.Range("D12").Value = "'40.80"
$ws.Range("E12").Value = "  +3.86%  "

# Row 13
$ws.Range("E13").Value = "  -0.53%  "

# Row 14
$ws.Range("D14").Value = "3.922.25"
$ws.Range("E14").Value = "  +4.50%  "

# Row 15
$ws.Range("D15").Value = "'8.34"
$ws.Range("E15").Value = "  -1.09%  "

# Row 16
$ws.Range("D16").Value = "'19.53"
$ws.Range("E16").Value = "  +2.86%  "

# Row 17
$ws.Range("D17").Value = "3.428.13"
$ws.Range("E17").Value = "  +5.75%  "

# Row 18
$ws.Range("D18").Value = "61.813.89"
$ws.Range("E18").Value = "  +8.74%  "

# Row 19
$ws.Range("E19").Value = "  -1.50%  "

# Row 20
$ws.Range("D20").Value = "'11.02"
$ws.Range("E20").Value = "  -0.37%  "

# Row 21
$ws.Range("D21").Value = "'0.0000115"
$ws.Range("E21").Value = "  +3.26%  "

# Row 22
$ws.Range("D22").Value = "'3.30"
$ws.Range("E22").Value = "  -0.32%  "

# Row 23
$ws.Range("D23").Value = "'12.54"
$ws.Range("E23").Value = "  -4.36%  "

# Row 24
$ws.Range("D24").Value = "'294.68"
$ws.Range("E24").Value = "  +1.16%  "

# Row 25
$ws.Range("D25").Value = "'74.51"
$ws.Range("E25").Value = "  +0.86%  "

# Row 26
$ws.Range("D26").Value = "'3.14"
$ws.Range("E26").Value = "  -0.71%  "

# Row 27
$ws.Range("D27").Value = "'30.64"
$ws.Range("E27").Value = "  +9.45%  "

# Row 28
$ws.Range("D28").Value = "'7.85"
$ws.Range("E28").Value = "  +9.73%  "

# Row 29
$ws.Range("D29").Value = "'0.174"
$ws.Range("E29").Value = "  +2.90%  "

# Row 30
$ws.Range("D30").Value = "'7.68"
$ws.Range("E30").Value = "  -0.19%  "

# Row 31
$ws.Range("D31").Value = "'4.24"
$ws.Range("E31").Value = "  -2.56%  "

# Row 32
$ws.Range("D32").Value = "'43.32"
$ws.Range("E32").Value = "  +8.60%  "

# Row 33
$ws.Range("E33").Value = "  +3.97%  "

# Row 34
$ws.Range("E34").Value = "  +2.17%  "

# Row 35
$ws.Range("D35").Value = "'1.00"

# Row 36
$ws.Range("E36").Value = "  +18.73%  "

# Row 37
$ws.Range("D37").Value = "'0.0480"
$ws.Range("E37").Value = "  -0.86%  "

# Row 38
$ws.Range("D38").Value = "'52.47"
$ws.Range("E38").Value = "  +2.03%  "

# Row 39
$ws.Range("B39").Value = "FirstDigitalUSD"
$ws.Range("C39").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D39").Value = "'0.998"
$ws.Range("E39").Value = "  -0.12%  "

# Row 40
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "'3.06"
$ws.Range("E40").Value = "  +4.77%  "

# Row 41
$ws.Range("D41").Value = "'3.45"
$ws.Range("E41").Value = "  -0.12%  "

# Row 42
$ws.Range("D42").Value = "'134.04"
$ws.Range("E42").Value = "  -2.41%  "

# Row 43
$ws.Range("E43").Value = "  -1.14%  "

# Row 44
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").Value = "'1.91"
$ws.Range("E44").Value = "  +0.93%  "

# Row 45
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").Value = "'0.285"
$ws.Range("E45").Value = "  +3.97%  "

# Row 46
$ws.Range("B46").Value = "Celestia"
$ws.Range("C46").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D46").Value = "'17.00"
$ws.Range("E46").Value = "  -0.09%  "

# Row 47
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "'3.86"
$ws.Range("E47").Value = "  -4.56%  "

# Row 48
$ws.Range("D48").Value = "'2.24"
$ws.Range("E48").Value = "  -2.60%  "

# Row 49
$ws.Range("D49").Value = "2.189.79"
$ws.Range("E49").Value = "  +0.99%  "

# Row 50
$ws.Range("D50").Value = "'21.13"
$ws.Range("E50").Value = "  -4.92%  "

# Row 51
$ws.Range("D51").Value = "3.727.06"
$ws.Range("E51").Value = "  +4.53%  "
